$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows for "Tasas de captación marginal bancaria 2021 - Diaria"
# Serie, Pesos a 1 mes, Pesos a 3 meses, Pesos a 6 meses, Pesos a 1 año, UF a 1 año
$data = @(
    @("15-09-2021", 1.61, 2.29, 2.85, 3.63, -0.33),
    @("16-09-2021", 1.61, 2.29, 2.85, 3.63, -0.31),
    @("20-09-2021", 1.65, 2.29, 2.89, 3.66, -0.26),
    @("21-09-2021", 1.75, 2.28, 2.91, 3.67, -0.22),
    @("22-09-2021", 1.78, 2.30, 2.95, 3.68, -0.21),
    @("23-09-2021", 1.66, 2.34, 2.98, 3.71, -0.21),
    @("24-09-2021", 1.66, 2.37, 3.03, 3.73, -0.18),
    @("27-09-2021", 1.60, 2.37, 3.02, 3.79, -0.21),
    @("28-09-2021", 1.71, 2.46, 3.06, 3.79, -0.20),
    @("29-09-2021", 1.73, 2.48, 3.08, 3.78, -0.16),
    @("30-09-2021", 1.73, 2.45, 3.12, 3.82, -0.16)
)

$startRow = 179
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
